$wb = $excel.ActiveWorkbook

# "展览" sheet (sheet1): update F2 and F5 (想去人数 counts)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 9396
$wsExhibit.Range("F5").Value = 516

# "全部类型" sheet (sheet4): update F2 and F5 (想去人数 counts)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 9396
$wsAll.Range("F5").Value = 516
